$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "Andrei Marica (?)" as trainer for "Know your language - OOP" row (row 8, column C)
$ws.Range("C8").Value = "Andrei Marica (?)"

# Add "Andrei Marica" as assistant for several session rows (column D)
$ws.Range("D9").Value = "Andrei Marica"
$ws.Range("D10").Value = "Andrei Marica"
$ws.Range("D17").Value = "Andrei Marica"
$ws.Range("D18").Value = "Andrei Marica"
$ws.Range("D19").Value = "Andrei Marica"
$ws.Range("D20").Value = "Andrei Marica"

# Update the active selection to C10, matching the author's final cursor position
$ws.Range("C10").Select()
